$d = $word.ActiveDocument

# --- Edit 1: Introduction paragraph -------------------------------------
# "...implications for social intelligence and collective cognition1."
# becomes
# "...implications for social intelligence, collective cognition, and
#  potential applications in engineering, artificial intelligence, and
#  robotics1."
$null = $d.Content.Find.Execute(
    "cations for social intelligence and collective cognition",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "cations for social intelligence, collective cognition, and potential applications in engineering, artificial intelligence, and robotics",
    2)

Write-Host "done edit 1"

# --- Edit 2: Intellectual Merit paragraph -------------------------------
# ", offering behavior researchers unprecedented experimental control
#  over predator-prey dynamics and the ability to test hypotheses about
#  the evolution of behavior in response to predation."
# becomes (word "behavior" removed after "offering"; several spans bolded)
$null = $d.Content.Find.Execute(
    ", offering behavior researchers",
    $true, $false, $false, $false, $false, $true, 1, $false,
    ", offering researchers",
    2)

$rngB1 = $d.Content
$null = $rngB1.Find.Execute("offering researchers unprecedented experimental control")
$rngB1.Bold = 1

$rngB2 = $d.Content
$null = $rngB2.Find.Execute("and the ability to test hypotheses about the evolution of behavior")
$rngB2.Bold = 1

Write-Host "done edit 2"

# --- Edit 3: Hypothesis I paragraph -------------------------------------
# Split the run "the predator than prey that move around randomly." into
# two runs ("the" / " predator than prey that move around randomly.")
# with no visible text change -- force a run break via a no-op bold
# toggle on the first word.
$rngSplit = $d.Content
$null = $rngSplit.Find.Execute("the predator than prey that move around randomly.")
$splitStart = $rngSplit.Start
$firstWord = $d.Range($splitStart, $splitStart + 3)
$firstWord.Bold = 1
$firstWord.Bold = 0

Write-Host "done edit 3"

# --- Edit 4: Broader Impacts paragraph ----------------------------------
# Remove " By doing so, this work will inform engineers, AI researchers,
# and roboticists on how to better design control algorithms for complex,
# distributed systems1" and replace with a new (partly bold) sentence.
$rngBI = $d.Content
$null = $rngBI.Find.Execute("By doing so, this work will inform engineers, AI researchers, and roboticists on how to better design control algorithms for complex, distributed systems1")
$rngBI.Text = ""
$rngBI.InsertAfter("The availability of such an interface should open up a new frontier in the study of the evolution of animal behavior and artificial intelligence")

$rngBIBold = $d.Content
$null = $rngBIBold.Find.Execute("The availability of such an interface should open up a new frontier")
$rngBIBold.Bold = 1

# "...through the NSF BEACON Center, local science fairs" -> "...through
# blog and video blog posts on my personal blog, science fairs"
$null = $d.Content.Find.Execute(
    "through the NSF BEACON Center, local science fairs",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "through blog and video blog posts on my personal blog, science fairs",
    2)

# "...and volunteering at the local museum." -> "...and volunteering at
# local museums."
$null = $d.Content.Find.Execute(
    "and volunteering at the local museum.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "and volunteering at local museums.",
    2)

Write-Host "done edit 4"

# --- Edit 5: move the _GoBack bookmark ----------------------------------
# It used to sit right after "stickleback's behavior" in the Intellectual
# Merit paragraph; it now sits right before "science fairs" in the
# Broader Impacts paragraph.
$d.Bookmarks("_GoBack").Delete()
$rngBookmark = $d.Content
$null = $rngBookmark.Find.Execute("science fairs")
$rngBookmark.Collapse(1)
$d.Bookmarks.Add("_GoBack", $rngBookmark)

Write-Host "done edit 5"
